$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''24.843.86'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '''1.716.47'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  +0.69%  '
$ws.Range('D5').Value = '''311.70'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '''1.005'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').Value = '''0.3777'
$ws.Range('E7').Value = '  +1.50%  '
$ws.Range('D8').Value = '''49.99'
$ws.Range('E8').Value = '  +3.69%  '
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').Value = '''1.194'
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +2.50%  '
$ws.Range('D12').Value = '''1.007'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = '''6.344'
$ws.Range('E13').Value = '  +4.17%  '
$ws.Range('D14').Value = '''20.86'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').Value = '''6.987'
$ws.Range('E15').Value = '  +3.56%  '
$ws.Range('D16').Value = '''1.718.90'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').Value = '''0.00001129'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').Value = '''1.006'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '''0.06681'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '''84.04'
$ws.Range('E20').Value = '  +3.18%  '
$ws.Range('D21').Value = '''17.30'
$ws.Range('E21').Value = '  +5.15%  '
$ws.Range('E22').Value = '  +4.18%  '
$ws.Range('D23').Value = '''13.43'
$ws.Range('E23').Value = '  +11.68%  '
$ws.Range('D24').Value = '''24.842.51'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').Value = '''2.450'
$ws.Range('E25').Value = '  +2.05%  '
$ws.Range('D26').Value = '''2.817'
$ws.Range('E26').Value = '  +5.56%  '
$ws.Range('D27').Value = '''20.49'
$ws.Range('E27').Value = '  +5.00%  '
$ws.Range('D28').Value = '''150.81'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''132.09'
$ws.Range('E29').Value = '  +4.05%  '
$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').Value = '''1.911.24'
$ws.Range('E30').Value = '  +2.71%  '
$ws.Range('D31').Value = '''1.185'
$ws.Range('E31').Value = '  +20.12%  '
$ws.Range('E32').Value = '  +6.42%  '
$ws.Range('D33').Value = '''4.234'
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '''0.08868'
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('D35').Value = '''13.71'
$ws.Range('E35').Value = '  +10.05%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.781'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +5.13%  '
$ws.Range('D38').Value = '''0.02433'
$ws.Range('E38').Value = '  +4.24%  '
$ws.Range('D39').Value = '''0.06533'
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('D40').Value = '''8.995'
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('D41').Value = '''0.2207'
$ws.Range('E41').Value = '  +4.63%  '
$ws.Range('D42').Value = '''1.279'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('D43').Value = '''0.6446'
$ws.Range('E43').Value = '  +4.97%  '
$ws.Range('D44').Value = '''1.006'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '''13.92'
$ws.Range('E45').Value = '  +6.00%  '
$ws.Range('D46').Value = '''0.6152'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('D47').Value = '''3.843'
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('D48').Value = '''2.149'
$ws.Range('E48').Value = '  +6.50%  '
$ws.Range('D49').Value = '''129.18'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = '''0.07281'
$ws.Range('E50').Value = '  +1.78%  '
$ws.Range('D51').Value = '''79.93'
$ws.Range('E51').Value = '  +4.38%  '
